$d = $word.ActiveDocument

# Locate the paragraph "Fundamentos em ecologia. Artmed." (last bibliography
# line). We want to delete everything from right after its paragraph mark
# through (and including) the paragraph mark of the copyright/footer
# paragraph, i.e. remove the blank paragraph, the "Ver no Jupiter..." line,
# and the "(c) 2020 ..." line that follow it.
$rStart = $d.Content
$foundStart = $rStart.Find.Execute("Fundamentos em ecologia. Artmed.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundStart) {
    throw "Could not find the 'Fundamentos em ecologia. Artmed.' paragraph"
}
$rStart.MoveEnd(1, 1) | Out-Null   # step over this paragraph's own mark
$startDel = $rStart.End

$rEnd = $d.Content
$foundEnd = $rEnd.Find.Execute("© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundEnd) {
    throw "Could not find the copyright/footer paragraph"
}
$rEnd.MoveEnd(1, 1) | Out-Null     # include this paragraph's own mark
$endDel = $rEnd.End

$d.Range($startDel, $endDel).Delete()
